$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "TestCases": mark existing RediffLogin-ish row off (B3 Y->N),
# add a new "RediffLogin" test case row (row 6).
# ---------------------------------------------------------------------
$wsTestCases = $wb.Worksheets.Item("TestCases")
$wsTestCases.Range("B3").Value = "N"
$wsTestCases.Range("A6").Value = "RediffLogin"
$wsTestCases.Range("B6").Value = "Y"

# ---------------------------------------------------------------------
# Sheet "Keywords": add a new keyword block header for RediffLogin.
# ---------------------------------------------------------------------
$wsKeywords = $wb.Worksheets.Item("Keywords")
$wsKeywords.Range("A20").Value = "RediffLogin"

# ---------------------------------------------------------------------
# Sheet "Data": add a new RediffLogin data section (header + column
# titles + one data row), reusing the formatting of the existing
# "TestD" section above it.
# ---------------------------------------------------------------------
$wsData = $wb.Worksheets.Item("Data")

# Section title (row 21) - copy format from the "TestD" section title (A16)
$wsData.Range("A16").Copy()
$wsData.Range("A21").PasteSpecial(-4122)
$wsData.Range("A21").Value = "RediffLogin"

# Column header row (row 22) - copy format from row 17 / row 13 headers
$wsData.Range("A17").Copy()
$wsData.Range("A22").PasteSpecial(-4122)
$wsData.Range("A22").Value = "Runmode"

$wsData.Range("B13").Copy()
$wsData.Range("B22").PasteSpecial(-4122)
$wsData.Range("B22").Value = "Browser"

$wsData.Range("C13").Copy()
$wsData.Range("C22").PasteSpecial(-4122)
$wsData.Range("C22").Value = "Username"

$wsData.Range("D13").Copy()
$wsData.Range("D22").PasteSpecial(-4122)
$wsData.Range("D22").Value = "Password"

# Data row (row 23)
$wsData.Range("A23").Value = "Y"

# ---------------------------------------------------------------------
# Restore selections per sheet, and leave "Data" as the active tab.
# ---------------------------------------------------------------------
[void]$wsTestCases.Activate()
[void]$wsTestCases.Range("A6").Select()

[void]$wsKeywords.Activate()
[void]$wsKeywords.Range("A20").Select()

[void]$wsData.Activate()
[void]$wsData.Range("A21").Select()
